$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching the style of existing
# header cells (bold/centered/bordered style used in B1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns (rows 2-27)
$data = @(
    @(2,9,9),
    @(3,6,7),
    @(4,4,5),
    @(5,7,8),
    @(6,3,3),
    @(7,5,6),
    @(8,9,9),
    @(9,7,8),
    @(10,8,8),
    @(11,7,7),
    @(12,8,8),
    @(13,5,5),
    @(14,1,2),
    @(15,7,7),
    @(16,6,6),
    @(17,8,9),
    @(18,5,5),
    @(19,9,9),
    @(20,5,8),
    @(21,8,8),
    @(22,5,6),
    @(23,5,6),
    @(24,10,10),
    @(25,9,9),
    @(26,7,7),
    @(27,6,6)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
